$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 9: change "Approved/Rejected" answer from Rejected to Approved,
# and clear the now-unused ReasonToReject ("Nil") cell.
$ws.Range("I9").Value = "Approved"
$ws.Range("J9").ClearContents()

# Update the sheet's selection to match the new review state
# (active cell I9, with I7:J7 also part of the selected set).
$ws.Range("I9:J9,I7:J7").Select()
